# rs data onboarding tracking
# Adds 8 new monitoring-well rows (13-20) to Sheet1: Piney 1/2/4/5 and
# P Port 1-4, each with Latitude (A), Longitude (B) and Name (C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row, Latitude, Longitude, Name
$newRows = @(
    @(13, 27.6358,             -82.5653,             "Piney 4"),
    @(14, 27.64742,             -82.56516,            "Piney 1"),
    @(15, 27.64948,             -82.57845,            "Piney 2"),
    @(16, 27.62852,             -82.59193,            "Piney 5"),
    @(17, 27.63488,             -82.56319,            "P Port 1"),
    @(18, 27.63105,             -82.55759,            "P Port 2"),
    @(19, 27.63141,             -82.54427,            "P Port 3"),
    @(20, 27.6454,              -82.53695,            "P Port 4")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Move the view roughly where the author left it (row 7 at top, D9 selected)
$ws.Range("D9").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
